$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of Argent (silver) price data appended below the existing table.
$newRow = 66

# Force the cells to be treated as plain text so values like "5,360" and
# "2025-05-06" are stored literally instead of being auto-converted into
# numbers/dates by Excel's input parsing.
$rng = $ws.Range("A$newRow`:J$newRow")
$rng.NumberFormat = "@"

$ws.Range("A$newRow").Value = "2025-05-06"
$ws.Range("B$newRow").Value = "38"
$ws.Range("C$newRow").Value = "37.3"
$ws.Range("D$newRow").Value = "1.02"
$ws.Range("E$newRow").Value = "0.273"
$ws.Range("F$newRow").Value = "0.09"
$ws.Range("G$newRow").Value = "5,360"
$ws.Range("H$newRow").Value = "8,025"
$ws.Range("I$newRow").Value = "8,075"
$ws.Range("J$newRow").Value = "7.2927"

# Restore the default (Normal) cell style so the new row matches the
# unstyled look of the rest of the worksheet.
$rng.Style = "Normal"
